# "Added the last assesments" - fill in the remaining blanks on the
# "Peer  and self assessment" sheet:
#   - Self-assessment block (row 2): replace the example action text in C2
#     with the actual action taken.
#   - Peer-assessment block (row 15, the "Self assesment" placeholder row):
#     grade the collaborator as "Good" and note the example action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

$ws.Range("C2").Value = "Support to others in github, Makeing sure that the git folder is almost clean so it is esay to find stuff"

$ws.Range("B15").Value = "Good"
$ws.Range("C15").Value = "Active in meetings, responding fast in discord"

$ws.Range("C17").Select() | Out-Null
